$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.498.17"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "3.337.39"
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "559.97"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").Value = "151.38"
$ws.Range("E6").Value = "  +4.18%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.340.09"
$ws.Range("E8").Value = "  +4.13%  "
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").Value = "7.39"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "3.913.66"
$ws.Range("E13").Value = "  +4.25%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "62.527.99"
$ws.Range("E17").Value = "  +4.07%  "
$ws.Range("D18").Value = "3.302.80"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "13.79"
$ws.Range("E20").Value = "  +4.42%  "
$ws.Range("D21").Value = "8.38"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").Value = "383.28"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "69.97"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "0.178"
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("D27").Value = "8.99"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "0.0₃0945"
$ws.Range("E29").Value = "  +5.10%  "
$ws.Range("D30").Value = "6.57"
$ws.Range("E30").Value = "  +6.49%  "
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "5.58"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "22.85"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +7.66%  "
$ws.Range("D35").Value = "6.72"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "159.36"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").Value = "  +8.31%  "
$ws.Range("E38").Value = "  +12.25%  "
$ws.Range("D39").Value = "26.81"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("D40").Value = "0.0737"
$ws.Range("E40").Value = "  +4.74%  "
$ws.Range("D41").Value = "2.794.13"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "0.0314"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "4.25"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "40.42"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.742"
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("E46").Value = "  +4.66%  "
$ws.Range("D47").Value = "3.380.28"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").Value = "21.92"
$ws.Range("E48").Value = "  +6.17%  "
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").Value = "287.47"
$ws.Range("E51").Value = "  +6.03%  "
